$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pork Sandwich recipe and advancement are now done ("Yes")
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "Yes"

# Move active selection to E2 (reflecting the latest edited cell)
$ws.Range("E2").Select()
